$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'65.413.73"
$ws.Range("E2").Value = "  -0.83%  "
$ws.Range("D3").Value = "'3.329.05"
$ws.Range("E3").Value = "  -4.17%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'575.97"
$ws.Range("E5").Value = "  -1.09%  "
$ws.Range("D6").Value = "'178.55"
$ws.Range("E6").Value = "  +3.07%  "
$ws.Range("E7").Value = "  +2.67%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").Value = "'3.327.31"
$ws.Range("E9").Value = "  -4.12%  "
$ws.Range("E10").Value = "  -1.16%  "
$ws.Range("E11").Value = "  +0.12%  "
$ws.Range("D12").Value = "'0.407"
$ws.Range("E12").Value = "  -0.71%  "
$ws.Range("D13").Value = "'3.909.50"
$ws.Range("E13").Value = "  -3.88%  "
$ws.Range("E14").Value = "  +0.39%  "
$ws.Range("D15").Value = "'28.67"
$ws.Range("E15").Value = "  -3.91%  "
$ws.Range("D16").Value = "'65.456.09"
$ws.Range("E16").Value = "  -0.90%  "
$ws.Range("E17").Value = "  -1.26%  "
$ws.Range("D18").Value = "'3.330.29"
$ws.Range("E18").Value = "  -3.90%  "
$ws.Range("D19").Value = "'5.73"
$ws.Range("E19").Value = "  -3.22%  "
$ws.Range("D20").Value = "'13.40"
$ws.Range("E20").Value = "  -3.25%  "
$ws.Range("D21").Value = "'362.30"
$ws.Range("E21").Value = "  -1.13%  "
$ws.Range("E22").Value = "  -3.89%  "
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("D24").Value = "'71.31"
$ws.Range("E24").Value = "  -1.12%  "
$ws.Range("E25").Value = "  -2.80%  "
$ws.Range("E26").Value = "  -1.94%  "
$ws.Range("D27").Value = "'9.54"
$ws.Range("E27").Value = "  -0.26%  "
$ws.Range("E28").Value = "  -0.37%  "
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("E30").Value = "  -0.97%  "
$ws.Range("D31").Value = "'5.60"
$ws.Range("E31").Value = "  -2.69%  "
$ws.Range("E32").Value = "  -0.04%  "
$ws.Range("E33").Value = "  -4.12%  "
$ws.Range("D34").Value = "'6.82"
$ws.Range("E34").Value = "  -4.00%  "
$ws.Range("D35").Value = "'1.21"
$ws.Range("E35").Value = "  -5.82%  "
$ws.Range("D36").Value = "'1.48"
$ws.Range("E36").Value = "  -2.65%  "
$ws.Range("D37").Value = "'160.78"
$ws.Range("E37").Value = "  +0.96%  "
$ws.Range("E38").Value = "  -5.20%  "
$ws.Range("D39").Value = "'27.31"
$ws.Range("E39").Value = "  -6.50%  "
$ws.Range("E40").Value = "  -0.82%  "
$ws.Range("E41").Value = "  +0.82%  "
$ws.Range("D42").Value = "'2.715.92"
$ws.Range("E42").Value = "  -3.14%  "
$ws.Range("D43").Value = "'6.24"
$ws.Range("E43").Value = "  -3.04%  "
$ws.Range("E44").Value = "  -3.63%  "
$ws.Range("D45").Value = "'336.86"
$ws.Range("E45").Value = "  +7.58%  "
$ws.Range("E46").Value = "  -2.20%  "
$ws.Range("D47").Value = "'39.72"
$ws.Range("E47").Value = "  -0.81%  "
$ws.Range("D48").Value = "'23.95"
$ws.Range("E48").Value = "  -0.83%  "
$ws.Range("D49").Value = "'0.0278"
$ws.Range("E49").Value = "  -3.28%  "
$ws.Range("E50").Value = "  +2.19%  "
$ws.Range("D51").Value = "'0.965"
$ws.Range("E51").Value = "  -0.78%  "
